{"js": "// The \"Error codes:\" paragraph currently reads (line breaks shown as \\v):\n//   Error codes:\n//   0 -> goed, geen error\n//   1 -> illegale zet\n//   ...\n//\n// The edit adds a new \"-1\" error-code line right before the existing\n// \"0\" line, and appends \"(move)\" to the existing \"0\" line, i.e. the\n// paragraph should become:\n//   Error codes:\n//   -1 -> goed, geen error (slag)\n//   0 -> goed, geen error (move)\n//   1 -> illegale zet\n//   ...\n\nconst body = context.document.body;\n\n// Step 1: insert a brand-new line \"-1 -> goed, geen error (slag)\" right\n// before the existing \"0 -> goed, geen error\" line (keeping the line\n// break that already precedes \"0 -> goed, geen error\").\nlet hits = body.search(\"0 -> goed, geen error\", { matchCase: true, matchWholeWord: false });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error('Could not find \"0 -> goed, geen error\" in the document.');\n}\n\nhits.items[0].insertText(\"-1 -> goed, geen error (slag)\\v\", \"Before\");\nawait context.sync();\n\n// Step 2: append \" (move)\" right after the (still intact) \"0 -> goed,\n// geen error\" text.\nhits = body.search(\"0 -> goed, geen error\", { matchCase: true, matchWholeWord: false });\nhits.load(\"items\");\nawait context.sync();\n\nhits.items[0].insertText(\" (move)\", \"After\");\nawait context.sync();\n", "ps1": "# The \"Error codes:\" paragraph currently reads (line breaks shown as `v):\n#   Error codes:\n#   0 -> goed, geen error\n#   1 -> illegale zet\n#   ...\n#\n# The edit adds a new \"-1\" error-code line right before the existing\n# \"0\" line, and appends \"(move)\" to the existing \"0\" line, i.e. the\n# paragraph should become:\n#   Error codes:\n#   -1 -> goed, geen error (slag)\n#   0 -> goed, geen error (move)\n#   1 -> illegale zet\n#   ...\n\n$d = $word.ActiveDocument\n\n# Step 1: insert a brand-new line \"-1 -> goed, geen error (slag)\" right\n# before the existing \"0 -> goed, geen error\" line (keeping the line\n# break that already precedes \"0 -> goed, geen error\").\n$range1 = $d.Content\n$find1 = $range1.Find\n$find1.Text = \"0 -> goed, geen error\"\n$find1.MatchCase = $true\n$found1 = $find1.Execute()\nif ($found1) {\n  $range1.InsertBefore(\"-1 -> goed, geen error (slag)`v\")\n}\n\n# Step 2: append \" (move)\" right after the (still intact) \"0 -> goed,\n# geen error\" text. Re-run Find since the document shifted.\n$range2 = $d.Content\n$find2 = $range2.Find\n$find2.Text = \"0 -> goed, geen error\"\n$find2.MatchCase = $true\n$found2 = $find2.Execute()\nif ($found2) {\n  $range2.Collapse(0)  # wdCollapseEnd\n  $range2.InsertAfter(\" (move)\")\n}\n"}
